$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price/percentage cells so values are stored as literal strings
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "E22", "D23", "E23", "D25", "E25", "E26", "E27", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "E47", "D48", "E48", "E49", "E50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '244.95'
$ws.Range('E2').Value = '-0.79%'
$ws.Range('D3').Value = '27.49'
$ws.Range('E3').Value = '5.23%'
$ws.Range('D4').Value = '5.116'
$ws.Range('E4').Value = '0.59%'
$ws.Range('D5').Value = '0.05681'
$ws.Range('E5').Value = '1.60%'
$ws.Range('E6').Value = '0.73%'
$ws.Range('D7').Value = '0.8197'
$ws.Range('E7').Value = '0.79%'
$ws.Range('D8').Value = '0.8517'
$ws.Range('E8').Value = '0.80%'
$ws.Range('B9').Value = 'MandalaExchangeToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D9').Value = '0.06951'
$ws.Range('E9').Value = '-0.48%'
$ws.Range('B10').Value = 'BitrueCoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D10').Value = '0.02879'
$ws.Range('E10').Value = '2.08%'
$ws.Range('B11').Value = 'BitMartToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D11').Value = '0.09396'
$ws.Range('E11').Value = '0.14%'
$ws.Range('B12').Value = 'BitForexToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D12').Value = '0.001526'
$ws.Range('E12').Value = '0.72%'
$ws.Range('B13').Value = 'CoinExToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D13').Value = '0.04025'
$ws.Range('E13').Value = '-13.52%'
$ws.Range('B14').Value = 'One'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D14').Value = '0.0005980'
$ws.Range('E14').Value = '-0.39%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = '0.006215'
$ws.Range('E15').Value = '0.19%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = '3.512'
$ws.Range('E16').Value = '-2.67%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = '3.009'
$ws.Range('E17').Value = '-0.35%'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').Value = '2.230'
$ws.Range('E18').Value = '8.49%'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').Value = '0.3192'
$ws.Range('E19').Value = '2.57%'
$ws.Range('B20').Value = 'WazirX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D20').Value = '0.1334'
$ws.Range('E20').Value = '0.22%'
$ws.Range('D21').Value = '0.03236'
$ws.Range('E21').Value = '1.82%'
$ws.Range('E22').Value = '-0.08%'
$ws.Range('D23').Value = '3.558'
$ws.Range('E23').Value = '-5.48%'
$ws.Range('D25').Value = '0.001216'
$ws.Range('E25').Value = '-2.25%'
$ws.Range('E26').Value = '-1.85%'
$ws.Range('E27').Value = '22.90%'
$ws.Range('E28').Value = '-27.48%'
$ws.Range('D40').Value = '0.03720'
$ws.Range('E40').Value = '1.61%'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = '0.1060'
$ws.Range('E41').Value = '-21.50%'
$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D42').Value = '0.003450'
$ws.Range('E42').Value = '-44.00%'
$ws.Range('D43').Value = '0.002380'
$ws.Range('E43').Value = '-6.98%'
$ws.Range('D44').Value = '0.009706'
$ws.Range('E44').Value = '20.53%'
$ws.Range('D45').Value = '0.00005099'
$ws.Range('E45').Value = '-5.40%'
$ws.Range('E46').Value = '-0.07%'
$ws.Range('E47').Value = '-30.39%'
$ws.Range('D48').Value = '0.002518'
$ws.Range('E48').Value = '4.73%'
$ws.Range('E49').Value = '-0.07%'
$ws.Range('E50').Value = '-0.07%'
